# Update manual testing documents (Refs #0).
#
# The "issues" sheet assigns two unscripted-testing tasks (rows 6 and 7,
# column C) away from their previous assignees over to "peterfpeterson".
# The dependent COUNTIF() totals on the "assignees" sheet recalculate
# automatically.  Finally, move the saved cell-selection/cursor position
# on the "issues" sheet.

$wb = $excel.ActiveWorkbook

$issues = $wb.Worksheets.Item("issues")

# Row 6 ("Unscripted Testing QECoverage") was assigned to rosswhitfield ->
# now peterfpeterson.
$issues.Range("C6").Value = "peterfpeterson"

# Row 7 ("Unscripted Testing ORNL HFIR diffraction & 4Circle") was assigned
# to wdzhou -> now peterfpeterson.
$issues.Range("C7").Value = "peterfpeterson"

# Move the active cell / selection that gets persisted with the sheet view.
$issues.Range("D17").Select()
